# Weekly update: insert a new price record as the first data row (row 448)
# for "Hortaliza, Feria Lagunitas de Puerto Montt - Cilantro", pushing all
# subsequent rows (448-549) down by one (to 449-550).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 448, shifting rows 448:549
# down to 449:550 (this also extends the sheet dimension to A1:R550).
$ws.Rows.Item(448).Insert()

# Populate the newly inserted row 448 with the new weekly record.
$ws.Range("A448").Value = 4
$ws.Range("B448").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C448").Value = "Los Lagos"
$ws.Range("D448").Value = 45244
$ws.Range("E448").Value = 10
$ws.Range("F448").Value = 100112040
$ws.Range("G448").Value = "Cilantro"
$ws.Range("H448").Value = "Sin especificar"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 160
$ws.Range("K448").Value = 9000
$ws.Range("L448").Value = 9000
$ws.Range("M448").Value = 9000
$ws.Range("N448").Value = "$/docena de atados (2 kilos)"
$ws.Range("O448").Value = "Región de La Araucanía"
$ws.Range("P448").Value = 4500
$ws.Range("Q448").Value = 2
$ws.Range("R448").Value = "Hortaliza"
